$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.443.87"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").Value = "2.513.97"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.143"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.81%  "
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "3.045.40"
$ws.Range("E13").Value = "  +3.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").Value = "68.237.13"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "2.505.78"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("D27").Value = "2.606.19"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").Value = "0.0₃0900"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "508.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").Value = "  +2.57%  "
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.26%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "150.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.96%  "
$ws.Range("E46").Value = "  +1.98%  "
$ws.Range("D47").Value = "0.0₆0259"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.518"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.579"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.87%  "
